$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 446.70587
$ws.Range("I39").Value = 397.45456
$ws.Range("K39").Value = 1192.36368
$ws.Range("M39").Value = -896.3636799999999
$ws.Range("H75").Value = 85000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
$ws.Range("H78").Value = 85000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
$ws.Range("H80").Value = 1309.25
$ws.Range("I80").Value = 1133.1666
$ws.Range("K80").Value = 3399.4998
$ws.Range("M80").Value = -2401.4998
$ws.Range("H83").Value = 1309.25
$ws.Range("I83").Value = 1133.1666
$ws.Range("K83").Value = 10198.4994
$ws.Range("M83").Value = -5206.499400000001
$ws.Range("H86").Value = 41278.6
$ws.Range("I86").Value = 63133
$ws.Range("K86").Value = 63133
$ws.Range("M86").Value = -62010
$ws.Range("H89").Value = 41278.6
$ws.Range("I89").Value = 63133
$ws.Range("K89").Value = 315665
$ws.Range("M89").Value = -310049
$ws.Range("H95").Value = 24492
$ws.Range("J95").Value = 24492
$ws.Range("L95").Value = 24492
$ws.Range("N95").Value = -29984

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2216.7856
$ws.Range("I61").Value = 1862.6364
$ws.Range("J61").Value = 3515.3333
$ws.Range("K61").Value = 1862.6364
$ws.Range("L61").Value = 3515.3333
$ws.Range("M61").Value = -1650.6364
$ws.Range("N61").Value = -3939.3333
$ws.Range("H74").Value = 1911
$ws.Range("I74").Value = 1886.1428
$ws.Range("J74").Value = 1998
$ws.Range("K74").Value = 1886.1428
$ws.Range("L74").Value = 1998
$ws.Range("M74").Value = -1012.1428
$ws.Range("N74").Value = -3746
$ws.Range("H77").Value = 1911
$ws.Range("I77").Value = 1886.1428
$ws.Range("J77").Value = 1998
$ws.Range("K77").Value = 9430.714
$ws.Range("L77").Value = 9990
$ws.Range("M77").Value = -5062.714
$ws.Range("N77").Value = -18726
$ws.Range("H92").Value = 49500
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 17000
$ws.Range("J96").Value = 17000
$ws.Range("L96").Value = 17000
$ws.Range("N96").Value = -22492
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H104").Value = 67498
$ws.Range("J104").Value = 67498
$ws.Range("L104").Value = 67498
$ws.Range("N104").Value = -74486
$ws.Range("H132").Value = 2516
$ws.Range("I132").Value = 2516
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7548
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5018
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2216.7856
$ws.Range("I136").Value = 1862.6364
$ws.Range("J136").Value = 3515.3333
$ws.Range("K136").Value = 5587.9092
$ws.Range("L136").Value = 10545.9999
$ws.Range("M136").Value = -3037.9092
$ws.Range("N136").Value = -15645.9999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2561.6667
$ws.Range("I105").Value = 2509.1667
$ws.Range("K105").Value = 2509.1667
$ws.Range("M105").Value = -762.1667000000002

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 710.8570999999999
$ws.Range("I105").Value = 710.8570999999999
$ws.Range("K105").Value = 710.8570999999999
$ws.Range("M105").Value = 1036.1429

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 561.13336
$ws.Range("I5").Value = 494
$ws.Range("J5").Value = 997.5
$ws.Range("K5").Value = 1482
$ws.Range("L5").Value = 2992.5
$ws.Range("M5").Value = -1370
$ws.Range("N5").Value = -3216.5
$ws.Range("H23").Value = 541.5
$ws.Range("I23").Value = 429
$ws.Range("J23").Value = 557.5714
$ws.Range("K23").Value = 1287
$ws.Range("L23").Value = 1672.7142
$ws.Range("M23").Value = -1052
$ws.Range("N23").Value = -2142.7142
$ws.Range("H80").Value = 5499.778
$ws.Range("I80").Value = 5250
$ws.Range("J80").Value = 5571.143
$ws.Range("K80").Value = 15750
$ws.Range("L80").Value = 16713.429
$ws.Range("M80").Value = -14814
$ws.Range("N80").Value = -18585.429
$ws.Range("H83").Value = 5499.778
$ws.Range("I83").Value = 5250
$ws.Range("J83").Value = 5571.143
$ws.Range("K83").Value = 47250
$ws.Range("L83").Value = 50140.287
$ws.Range("M83").Value = -42570
$ws.Range("N83").Value = -59500.287
$ws.Range("H135").Value = 561.13336
$ws.Range("I135").Value = 494
$ws.Range("J135").Value = 997.5
$ws.Range("K135").Value = 4446
$ws.Range("L135").Value = 8977.5
$ws.Range("M135").Value = -1911
$ws.Range("N135").Value = -14047.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 36000
$ws.Range("I20").Value = 36000
$ws.Range("J20").Value = 36000
$ws.Range("K20").Value = 36000
$ws.Range("L20").Value = 36000
$ws.Range("M20").Value = -35755
$ws.Range("N20").Value = -36490
$ws.Range("H24").Value = 12230.77
$ws.Range("J24").Value = 12230.77
$ws.Range("L24").Value = 12230.77
$ws.Range("N24").Value = -12576.77
$ws.Range("H97").Value = 1689.5
$ws.Range("J97").Value = 2943
$ws.Range("L97").Value = 2943
$ws.Range("N97").Value = -3935
$ws.Range("H101").Value = 57885.668
$ws.Range("J101").Value = 57885.668
$ws.Range("L101").Value = 57885.668
$ws.Range("N101").Value = -64375.668

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5543.909
$ws.Range("I7").Value = 5098.3
$ws.Range("K7").Value = 5098.3
$ws.Range("M7").Value = -4986.3
$ws.Range("H93").Value = 2256.4285
$ws.Range("I93").Value = 1965.8334
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1965.8334
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -717.8334
$ws.Range("N93").Value = -6496
$ws.Range("H105").Value = 39853.332
$ws.Range("J105").Value = 39853.332
$ws.Range("L105").Value = 39853.332
$ws.Range("N105").Value = -46841.332
$ws.Range("H126").Value = 5543.909
$ws.Range("I126").Value = 5098.3
$ws.Range("K126").Value = 15294.9
$ws.Range("M126").Value = -12824.9

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1504.7142
$ws.Range("I132").Value = 1421.3334
$ws.Range("J132").Value = 2005
$ws.Range("K132").Value = 4264.0002
$ws.Range("L132").Value = 6015
$ws.Range("M132").Value = -1734.0002
$ws.Range("N132").Value = -11075
$ws.Range("H136").Value = 3303.3462
$ws.Range("J136").Value = 3448.2222
$ws.Range("L136").Value = 10344.6666
$ws.Range("N136").Value = -15444.6666

Write-Host "Applied all changes"